$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value2 = 1.02 ; $ws.Cells.Item(2, 3).Value2 = 1.028973255428667 ; $ws.Cells.Item(2, 4).Value2 = 1.037059800123672 ; $ws.Cells.Item(2, 5).Value2 = 1.038560239073564 ; $ws.Cells.Item(2, 6).Value2 = 1.048804984141571 ; $ws.Cells.Item(2, 9).Value2 = 1.032172886504102 ; $ws.Cells.Item(2, 10).Value2 = 1.034122606098828 ; $ws.Cells.Item(2, 11).Value2 = 1.039851764756822 ; $ws.Cells.Item(2, 12).Value2 = 1.041347924427985 ; $ws.Cells.Item(2, 13).Value2 = 1.051563806205685 ; $ws.Cells.Item(2, 14).Value2 = 1.015373910130145
$ws.Cells.Item(3, 2).Value2 = 1.02 ; $ws.Cells.Item(3, 3).Value2 = 1.029812177506045 ; $ws.Cells.Item(3, 4).Value2 = 1.037685523650234 ; $ws.Cells.Item(3, 5).Value2 = 1.039297506843562 ; $ws.Cells.Item(3, 6).Value2 = 1.049605956652344 ; $ws.Cells.Item(3, 9).Value2 = 1.032269299084244 ; $ws.Cells.Item(3, 10).Value2 = 1.034602929769353 ; $ws.Cells.Item(3, 11).Value2 = 1.04028770566879 ; $ws.Cells.Item(3, 12).Value2 = 1.041895422868453 ; $ws.Cells.Item(3, 13).Value2 = 1.05217692659084 ; $ws.Cells.Item(3, 14).Value2 = 1.015534661474798
$ws.Cells.Item(4, 2).Value2 = 1.02 ; $ws.Cells.Item(4, 3).Value2 = 1.030355786119329 ; $ws.Cells.Item(4, 4).Value2 = 1.038090882950311 ; $ws.Cells.Item(4, 5).Value2 = 1.03977561991438 ; $ws.Cells.Item(4, 6).Value2 = 1.050125306162631 ; $ws.Cells.Item(4, 9).Value2 = 1.032330479878947 ; $ws.Cells.Item(4, 10).Value2 = 1.034913848438617 ; $ws.Cells.Item(4, 11).Value2 = 1.040569571493908 ; $ws.Cells.Item(4, 12).Value2 = 1.042250073805092 ; $ws.Cells.Item(4, 13).Value2 = 1.05257406652184 ; $ws.Cells.Item(4, 14).Value2 = 1.015638671310203
$ws.Cells.Item(5, 2).Value2 = 1.02 ; $ws.Cells.Item(5, 3).Value2 = 1.030584501767522 ; $ws.Cells.Item(5, 4).Value2 = 1.038261407625476 ; $ws.Cells.Item(5, 5).Value2 = 1.039976868429146 ; $ws.Cells.Item(5, 6).Value2 = 1.050343893968554 ; $ws.Cells.Item(5, 9).Value2 = 1.032355911206709 ; $ws.Cells.Item(5, 10).Value2 = 1.035044585332057 ; $ws.Cells.Item(5, 11).Value2 = 1.040688014474715 ; $ws.Cells.Item(5, 12).Value2 = 1.042399259175252 ; $ws.Cells.Item(5, 13).Value2 = 1.052741120564155 ; $ws.Cells.Item(5, 14).Value2 = 1.015682394862004
$ws.Cells.Item(6, 2).Value2 = 1.02 ; $ws.Cells.Item(6, 3).Value2 = 1.030622914787837 ; $ws.Cells.Item(6, 4).Value2 = 1.038290045958315 ; $ws.Cells.Item(6, 5).Value2 = 1.040010673530576 ; $ws.Cells.Item(6, 6).Value2 = 1.050380610613689 ; $ws.Cells.Item(6, 9).Value2 = 1.032360164266608 ; $ws.Cells.Item(6, 10).Value2 = 1.035066538157065 ; $ws.Cells.Item(6, 11).Value2 = 1.040707898398395 ; $ws.Cells.Item(6, 12).Value2 = 1.042424313295162 ; $ws.Cells.Item(6, 13).Value2 = 1.052769175275869 ; $ws.Cells.Item(6, 14).Value2 = 1.01568973609621
$ws.Cells.Item(7, 2).Value2 = 1.02 ; $ws.Cells.Item(7, 3).Value2 = 1.030358841513742 ; $ws.Cells.Item(7, 4).Value2 = 1.038093161073189 ; $ws.Cells.Item(7, 5).Value2 = 1.039778308027882 ; $ws.Cells.Item(7, 6).Value2 = 1.050128225951056 ; $ws.Cells.Item(7, 9).Value2 = 1.032330820830343 ; $ws.Cells.Item(7, 10).Value2 = 1.034915595248237 ; $ws.Cells.Item(7, 11).Value2 = 1.040571154346497 ; $ws.Cells.Item(7, 12).Value2 = 1.042252066875143 ; $ws.Cells.Item(7, 13).Value2 = 1.052576298329482 ; $ws.Cells.Item(7, 14).Value2 = 1.015639255555646
$ws.Cells.Item(8, 2).Value2 = 1.02 ; $ws.Cells.Item(8, 3).Value2 = 1.02925661260829 ; $ws.Cells.Item(8, 4).Value2 = 1.037271166943697 ; $ws.Cells.Item(8, 5).Value2 = 1.038809183267124 ; $ws.Cells.Item(8, 6).Value2 = 1.049075454410499 ; $ws.Cells.Item(8, 9).Value2 = 1.032205718586585 ; $ws.Cells.Item(8, 10).Value2 = 1.034284908662147 ; $ws.Cells.Item(8, 11).Value2 = 1.039999137115286 ; $ws.Cells.Item(8, 12).Value2 = 1.041532873897974 ; $ws.Cells.Item(8, 13).Value2 = 1.051770927156425 ; $ws.Cells.Item(8, 14).Value2 = 1.015428237956143
$ws.Cells.Item(9, 2).Value2 = 1.02 ; $ws.Cells.Item(9, 3).Value2 = 1.027320318097807 ; $ws.Cells.Item(9, 4).Value2 = 1.035826426490617 ; $ws.Cells.Item(9, 5).Value2 = 1.037109601567647 ; $ws.Cells.Item(9, 6).Value2 = 1.047228599550636 ; $ws.Cells.Item(9, 9).Value2 = 1.031976077546006 ; $ws.Cells.Item(9, 10).Value2 = 1.033174515480733 ; $ws.Cells.Item(9, 11).Value2 = 1.038989570591399 ; $ws.Cells.Item(9, 12).Value2 = 1.040268563925969 ; $ws.Cells.Item(9, 13).Value2 = 1.050354978072565 ; $ws.Cells.Item(9, 14).Value2 = 1.01505636617351
$ws.Cells.Item(10, 2).Value2 = 1.02 ; $ws.Cells.Item(10, 3).Value2 = 1.026033575071767 ; $ws.Cells.Item(10, 4).Value2 = 1.034865881489415 ; $ws.Cells.Item(10, 5).Value2 = 1.035982136896064 ; $ws.Cells.Item(10, 6).Value2 = 1.046003044544694 ; $ws.Cells.Item(10, 9).Value2 = 1.031816843019978 ; $ws.Cells.Item(10, 10).Value2 = 1.032434977961999 ; $ws.Cells.Item(10, 11).Value2 = 1.03831553703536 ; $ws.Cells.Item(10, 12).Value2 = 1.039427799294434 ; $ws.Cells.Item(10, 13).Value2 = 1.049413280041281 ; $ws.Cells.Item(10, 14).Value2 = 1.014808460046273
$ws.Cells.Item(11, 2).Value2 = 1.02 ; $ws.Cells.Item(11, 3).Value2 = 1.025477400187863 ; $ws.Cells.Item(11, 4).Value2 = 1.034450599647613 ; $ws.Cells.Item(11, 5).Value2 = 1.035495282842366 ; $ws.Cells.Item(11, 6).Value2 = 1.045473739956051 ; $ws.Cells.Item(11, 9).Value2 = 1.031746445054455 ; $ws.Cells.Item(11, 10).Value2 = 1.032114939049607 ; $ws.Cells.Item(11, 11).Value2 = 1.038023456607723 ; $ws.Cells.Item(11, 12).Value2 = 1.03906425847724 ; $ws.Cells.Item(11, 13).Value2 = 1.04900607446689 ; $ws.Cells.Item(11, 14).Value2 = 1.014701122242843
$ws.Cells.Item(12, 2).Value2 = 1.02 ; $ws.Cells.Item(12, 3).Value2 = 1.025270963070184 ; $ws.Cells.Item(12, 4).Value2 = 1.034296443691712 ; $ws.Cells.Item(12, 5).Value2 = 1.035314647747176 ; $ws.Cells.Item(12, 6).Value2 = 1.045277340484231 ; $ws.Cells.Item(12, 9).Value2 = 1.031720079208974 ; $ws.Cells.Item(12, 10).Value2 = 1.031996091632822 ; $ws.Cells.Item(12, 11).Value2 = 1.037914933428064 ; $ws.Cells.Item(12, 12).Value2 = 1.038929302299093 ; $ws.Cells.Item(12, 13).Value2 = 1.048854905495431 ; $ws.Cells.Item(12, 14).Value2 = 1.014661253791433
$ws.Cells.Item(13, 2).Value2 = 1.02 ; $ws.Cells.Item(13, 3).Value2 = 1.025315237671428 ; $ws.Cells.Item(13, 4).Value2 = 1.034329506191473 ; $ws.Cells.Item(13, 5).Value2 = 1.035353385316436 ; $ws.Cells.Item(13, 6).Value2 = 1.045319459398523 ; $ws.Cells.Item(13, 9).Value2 = 1.031725744576666 ; $ws.Cells.Item(13, 10).Value2 = 1.032021583459742 ; $ws.Cells.Item(13, 11).Value2 = 1.037938213420589 ; $ws.Cells.Item(13, 12).Value2 = 1.038958247254198 ; $ws.Cells.Item(13, 13).Value2 = 1.048887327863734 ; $ws.Cells.Item(13, 14).Value2 = 1.01466980562958
$ws.Cells.Item(14, 2).Value2 = 1.02 ; $ws.Cells.Item(14, 3).Value2 = 1.025460332940815 ; $ws.Cells.Item(14, 4).Value2 = 1.034437855053275 ; $ws.Cells.Item(14, 5).Value2 = 1.035480347312227 ; $ws.Cells.Item(14, 6).Value2 = 1.045457501262439 ; $ws.Cells.Item(14, 9).Value2 = 1.031744270066256 ; $ws.Cells.Item(14, 10).Value2 = 1.032105114479656 ; $ws.Cells.Item(14, 11).Value2 = 1.038014486682036 ; $ws.Cells.Item(14, 12).Value2 = 1.039053101335403 ; $ws.Cells.Item(14, 13).Value2 = 1.048993577041891 ; $ws.Cells.Item(14, 14).Value2 = 1.01469782666888
$ws.Cells.Item(15, 2).Value2 = 1.02 ; $ws.Cells.Item(15, 3).Value2 = 1.025549750986423 ; $ws.Cells.Item(15, 4).Value2 = 1.034504625436169 ; $ws.Cells.Item(15, 5).Value2 = 1.035558599903929 ; $ws.Cells.Item(15, 6).Value2 = 1.045542581013766 ; $ws.Cells.Item(15, 9).Value2 = 1.031755655510346 ; $ws.Cells.Item(15, 10).Value2 = 1.032156584619821 ; $ws.Cells.Item(15, 11).Value2 = 1.038061477024993 ; $ws.Cells.Item(15, 12).Value2 = 1.039111554589701 ; $ws.Cells.Item(15, 13).Value2 = 1.049059052025212 ; $ws.Cells.Item(15, 14).Value2 = 1.014715091581627
$ws.Cells.Item(16, 2).Value2 = 1.02 ; $ws.Cells.Item(16, 3).Value2 = 1.026070507645747 ; $ws.Cells.Item(16, 4).Value2 = 1.034893456041334 ; $ws.Cells.Item(16, 5).Value2 = 1.036014476318519 ; $ws.Cells.Item(16, 6).Value2 = 1.046038201771727 ; $ws.Cells.Item(16, 9).Value2 = 1.031821484644925 ; $ws.Cells.Item(16, 10).Value2 = 1.032456221921483 ; $ws.Cells.Item(16, 11).Value2 = 1.038334916953001 ; $ws.Cells.Item(16, 12).Value2 = 1.039451937303089 ; $ws.Cells.Item(16, 13).Value2 = 1.049440316799225 ; $ws.Cells.Item(16, 14).Value2 = 1.014815583897585
$ws.Cells.Item(17, 2).Value2 = 1.02 ; $ws.Cells.Item(17, 3).Value2 = 1.026397431726149 ; $ws.Cells.Item(17, 4).Value2 = 1.035137532085121 ; $ws.Cells.Item(17, 5).Value2 = 1.036300797227034 ; $ws.Cells.Item(17, 6).Value2 = 1.046349459582927 ; $ws.Cells.Item(17, 9).Value2 = 1.031862390243035 ; $ws.Cells.Item(17, 10).Value2 = 1.032644227182224 ; $ws.Cells.Item(17, 11).Value2 = 1.03850638087793 ; $ws.Cells.Item(17, 12).Value2 = 1.03966558959555 ; $ws.Cells.Item(17, 13).Value2 = 1.049679624212776 ; $ws.Cells.Item(17, 14).Value2 = 1.014878622355316
$ws.Cells.Item(18, 2).Value2 = 1.02 ; $ws.Cells.Item(18, 3).Value2 = 1.026588216735587 ; $ws.Cells.Item(18, 4).Value2 = 1.035279959176326 ; $ws.Cells.Item(18, 5).Value2 = 1.036467933063448 ; $ws.Cells.Item(18, 6).Value2 = 1.046531142923606 ; $ws.Cells.Item(18, 9).Value2 = 1.03188610999356 ; $ws.Cells.Item(18, 10).Value2 = 1.032753905295905 ; $ws.Cells.Item(18, 11).Value2 = 1.03860637161242 ; $ws.Cells.Item(18, 12).Value2 = 1.039790259026454 ; $ws.Cells.Item(18, 13).Value2 = 1.049819261693186 ; $ws.Cells.Item(18, 14).Value2 = 1.014915392280424
$ws.Cells.Item(19, 2).Value2 = 1.02 ; $ws.Cells.Item(19, 3).Value2 = 1.026653285676506 ; $ws.Cells.Item(19, 4).Value2 = 1.035328533532665 ; $ws.Cells.Item(19, 5).Value2 = 1.036524944001523 ; $ws.Cells.Item(19, 6).Value2 = 1.046593114558446 ; $ws.Cells.Item(19, 9).Value2 = 1.031894174083135 ; $ws.Cells.Item(19, 10).Value2 = 1.032791305679811 ; $ws.Cells.Item(19, 11).Value2 = 1.038640462216177 ; $ws.Cells.Item(19, 12).Value2 = 1.039832776456569 ; $ws.Cells.Item(19, 13).Value2 = 1.049866883487899 ; $ws.Cells.Item(19, 14).Value2 = 1.014927929971007
$ws.Cells.Item(20, 2).Value2 = 1.02 ; $ws.Cells.Item(20, 3).Value2 = 1.026362345952374 ; $ws.Cells.Item(20, 4).Value2 = 1.035111338642984 ; $ws.Cells.Item(20, 5).Value2 = 1.036270064273783 ; $ws.Cells.Item(20, 6).Value2 = 1.046316050904174 ; $ws.Cells.Item(20, 9).Value2 = 1.031858015913905 ; $ws.Cells.Item(20, 10).Value2 = 1.032624054143213 ; $ws.Cells.Item(20, 11).Value2 = 1.038487986611188 ; $ws.Cells.Item(20, 12).Value2 = 1.039642661578289 ; $ws.Cells.Item(20, 13).Value2 = 1.049653943243852 ; $ws.Cells.Item(20, 14).Value2 = 1.014871858853535
$ws.Cells.Item(21, 2).Value2 = 1.02 ; $ws.Cells.Item(21, 3).Value2 = 1.025417601817952 ; $ws.Cells.Item(21, 4).Value2 = 1.03440594628698 ; $ws.Cells.Item(21, 5).Value2 = 1.035442954516418 ; $ws.Cells.Item(21, 6).Value2 = 1.045416845619857 ; $ws.Cells.Item(21, 9).Value2 = 1.031738820753417 ; $ws.Cells.Item(21, 10).Value2 = 1.032080515856728 ; $ws.Cells.Item(21, 11).Value2 = 1.037992026965498 ; $ws.Cells.Item(21, 12).Value2 = 1.039025166977378 ; $ws.Cells.Item(21, 13).Value2 = 1.04896228693775 ; $ws.Cells.Item(21, 14).Value2 = 1.014689575123295
$ws.Cells.Item(22, 2).Value2 = 1.02 ; $ws.Cells.Item(22, 3).Value2 = 1.024824478184719 ; $ws.Cells.Item(22, 4).Value2 = 1.033963007312013 ; $ws.Cells.Item(22, 5).Value2 = 1.034924100089419 ; $ws.Cells.Item(22, 6).Value2 = 1.044852683021822 ; $ws.Cells.Item(22, 9).Value2 = 1.031662623414306 ; $ws.Cells.Item(22, 10).Value2 = 1.031738941881559 ; $ws.Cells.Item(22, 11).Value2 = 1.037680015904177 ; $ws.Cells.Item(22, 12).Value2 = 1.03863738185403 ; $ws.Cells.Item(22, 13).Value2 = 1.048527909628801 ; $ws.Cells.Item(22, 14).Value2 = 1.014574975546193
$ws.Cells.Item(23, 2).Value2 = 1.02 ; $ws.Cells.Item(23, 3).Value2 = 1.025138820763799 ; $ws.Cells.Item(23, 4).Value2 = 1.034197763028115 ; $ws.Cells.Item(23, 5).Value2 = 1.035199041896148 ; $ws.Cells.Item(23, 6).Value2 = 1.045151641491625 ; $ws.Cells.Item(23, 9).Value2 = 1.03170313575623 ; $ws.Cells.Item(23, 10).Value2 = 1.031920000121454 ; $ws.Cells.Item(23, 11).Value2 = 1.037845435590696 ; $ws.Cells.Item(23, 12).Value2 = 1.038842910124439 ; $ws.Cells.Item(23, 13).Value2 = 1.048758133791039 ; $ws.Cells.Item(23, 14).Value2 = 1.014635725885991
$ws.Cells.Item(24, 2).Value2 = 1.02 ; $ws.Cells.Item(24, 3).Value2 = 1.026378199407738 ; $ws.Cells.Item(24, 4).Value2 = 1.03512317414173 ; $ws.Cells.Item(24, 5).Value2 = 1.036283950772004 ; $ws.Cells.Item(24, 6).Value2 = 1.0463311464407 ; $ws.Cells.Item(24, 9).Value2 = 1.031859992916971 ; $ws.Cells.Item(24, 10).Value2 = 1.032633169416301 ; $ws.Cells.Item(24, 11).Value2 = 1.03849629825447 ; $ws.Cells.Item(24, 12).Value2 = 1.039653021609116 ; $ws.Cells.Item(24, 13).Value2 = 1.049665547203165 ; $ws.Cells.Item(24, 14).Value2 = 1.014874914986887
$ws.Cells.Item(25, 2).Value2 = 1.02 ; $ws.Cells.Item(25, 3).Value2 = 1.027820178066316 ; $ws.Cells.Item(25, 4).Value2 = 1.036199473998851 ; $ws.Cells.Item(25, 5).Value2 = 1.037548007491301 ; $ws.Cells.Item(25, 6).Value2 = 1.0477050635992 ; $ws.Cells.Item(25, 9).Value2 = 1.032036530615387 ; $ws.Cells.Item(25, 10).Value2 = 1.033461456951137 ; $ws.Cells.Item(25, 11).Value2 = 1.039250748378228 ; $ws.Cells.Item(25, 12).Value2 = 1.040595053398477 ; $ws.Cells.Item(25, 13).Value2 = 1.050720642745478 ; $ws.Cells.Item(25, 14).Value2 = 1.015152504642058
